# 파워레벨 테이블에 requiredLimitBreakGold|Int 컬럼 추가
# Add a new "requiredLimitBreakGold|Int" column to the PowerLevelTable sheet,
# inserted right after the existing "requiredLimitBreak|Int" column (column D),
# pushing the old E:G columns (requiredPowerPoint / requiredAccumulatedPowerPoint /
# requiredGold) one position to the right (F:H).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PowerLevelTable")

# Insert a brand-new column at E; this shifts old columns E:G -> F:H and leaves
# the new column E empty/unformatted, exactly like the target workbook.
$ws.Columns.Item(5).Insert()

# Header for the newly inserted column (this text becomes a new shared string).
$ws.Range("E1").Value = "requiredLimitBreakGold|Int"

# The new column stores the gold cost required to perform a "limit break" at a
# given power level. Almost every row is 0, except the three rows where a new
# limit-break tier unlocks (power levels 11, 14 and 17 -> sheet rows 12, 15, 18),
# which carry the actual gold cost for that tier.
$limitBreakGoldByRow = @{
    12 = 270000
    15 = 465000
    18 = 795000
}

for ($row = 2; $row -le 21; $row++) {
    if ($limitBreakGoldByRow.ContainsKey($row)) {
        $ws.Cells.Item($row, 5).Value = $limitBreakGoldByRow[$row]
    } else {
        $ws.Cells.Item($row, 5).Value = 0
    }
}

# Match the author's final selection/active cell in the sheet.
$ws.Range("F1").Select()
